$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# 1. JavaMethodService.java:163 -> 162 and AbstractService.java:136 -> 135
$d.Content.Find.Execute("JavaMethodService.java:163", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "JavaMethodService.java:162", $wdReplaceOne) | Out-Null
$d.Content.Find.Execute("AbstractService.java:136", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "AbstractService.java:135", $wdReplaceOne) | Out-Null

# 2. EvaluationServices.java:168 -> 172
$d.Content.Find.Execute("EvaluationServices.java:168", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "EvaluationServices.java:172", $wdReplaceOne) | Out-Null

# 3. AstEvaluator.java:189 -> 186 (caseCall line)
$d.Content.Find.Execute("caseCall(AstEvaluator.java:189)", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "caseCall(AstEvaluator.java:186)", $wdReplaceOne) | Out-Null

# 4. AstSwitch.java:118 -> 119
$d.Content.Find.Execute("AstSwitch.java:118", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "AstSwitch.java:119", $wdReplaceOne) | Out-Null

# 5. AstEvaluator.java:112 -> 109 (eval line)
$d.Content.Find.Execute("eval(AstEvaluator.java:112)", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "eval(AstEvaluator.java:109)", $wdReplaceOne) | Out-Null

# 6. GeneratedMethodAccessor74 -> GeneratedMethodAccessor73
$d.Content.Find.Execute("GeneratedMethodAccessor74", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "GeneratedMethodAccessor73", $wdReplaceOne) | Out-Null

# 7. Replace the large trailing block of the stack trace (from JUnit4Provider.execute down to
#    equinox launcher Main.main) with the new shorter jdt junit runner block.
$oldBlock = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)"

$newBlock = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

$d.Content.Find.Execute($oldBlock, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newBlock, $wdReplaceOne) | Out-Null

Write-Host "Done"
